$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SNMP")

$ws.Range("F2").Value = "10.30.4.77"
$ws.Range("H2").Value = "rootpw"
$ws.Range("D2").Value = "DES|AES128"

$ws.Activate()
$ws.Range("E9").Select()
